$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -8.103999999999999
$ws.Range("D21").Value = -8.227
$ws.Range("D23").Value = -8.135
$ws.Range("D25").Value = -7.842999999999999
$ws.Range("D53").Value = -7.947999999999999
$ws.Range("D57").Value = -7.935
$ws.Range("D59").Value = -8.191000000000001
$ws.Range("D69").Value = -7.640000000000001
$ws.Range("D79").Value = -7.867
$ws.Range("D83").Value = -8.125
$ws.Range("D93").Value = -6.976999999999999
